$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 200.22223
$ws.Range("I5").Value = 200.25
$ws.Range("K5").Value = 200.25
$ws.Range("M5").Value = -85.25
$ws.Range("H17").Value = 377268.7
$ws.Range("J17").Value = 420661.22
$ws.Range("N17").Value = -1262319.66
$ws.Range("L17").Value = 1261983.66
$ws.Range("N57").Value = -253838
$ws.Range("L57").Value = 252840
$ws.Range("H57").Value = 84280
$ws.Range("J57").Value = 84280
$ws.Range("H137").Value = 4562.5273
$ws.Range("I137").Value = 3759.1052
$ws.Range("K137").Value = 11277.3156
$ws.Range("M137").Value = -8727.3156
$ws.Range("H138").Value = 1888.1731
$ws.Range("I138").Value = 772.4
$ws.Range("K138").Value = 2317.2
$ws.Range("J138").Value = 3409.682
$ws.Range("M138").Value = 2822.8
$ws.Range("N138").Value = -20509.046
$ws.Range("L138").Value = 10229.046
$ws.Range("N140").Value = -180796.25
$ws.Range("L140").Value = 170436.25
$ws.Range("H140").Value = 170436.25
$ws.Range("J140").Value = 170436.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 4661.9272
$ws.Range("N32").Value = -158698
$ws.Range("M32").Value = -4374.9272
$ws.Range("L32").Value = 158124
$ws.Range("H32").Value = 24149.174
$ws.Range("I32").Value = 4661.9272
$ws.Range("J32").Value = 158124
$ws.Range("M61").Value = -798
$ws.Range("H61").Value = 1151.4286
$ws.Range("I61").Value = 1010
$ws.Range("K61").Value = 1010
$ws.Range("H74").Value = 1579.3939
$ws.Range("J74").Value = 3180.6667
$ws.Range("N74").Value = -4928.6667
$ws.Range("L74").Value = 3180.6667
$ws.Range("J77").Value = 3180.6667
$ws.Range("N77").Value = -24639.3335
$ws.Range("L77").Value = 15903.3335
$ws.Range("H77").Value = 1579.3939
$ws.Range("H132").Value = 1675.5714
$ws.Range("I132").Value = 1401.0476
$ws.Range("K132").Value = 4203.142800000001
$ws.Range("J132").Value = 2499.1428
$ws.Range("N132").Value = -12557.4284
$ws.Range("L132").Value = 7497.428400000001
$ws.Range("M132").Value = -1673.142800000001
$ws.Range("M136").Value = -480
$ws.Range("H136").Value = 1151.4286
$ws.Range("I136").Value = 1010
$ws.Range("K136").Value = 3030
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M134").Value = -731.5715999999998
$ws.Range("H134").Value = 1062.9333
$ws.Range("I134").Value = 1088.8572
$ws.Range("K134").Value = 3266.5716
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I58").Value = 1609.6666
$ws.Range("K58").Value = 1609.6666
$ws.Range("J58").Value = 2702.6
$ws.Range("N58").Value = -3108.6
$ws.Range("L58").Value = 2702.6
$ws.Range("M58").Value = -1406.6666
$ws.Range("H58").Value = 1798.1034
$ws.Range("H105").Value = 2499
$ws.Range("I105").Value = 2398.8
$ws.Range("K105").Value = 2398.8
$ws.Range("J105").Value = 3000
$ws.Range("N105").Value = -6494
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -651.8000000000002
$ws.Range("H132").Value = 3772.1316
$ws.Range("I132").Value = 3661.4644
$ws.Range("K132").Value = 10984.3932
$ws.Range("J132").Value = 4082
$ws.Range("N132").Value = -17306
$ws.Range("L132").Value = 12246
$ws.Range("M132").Value = -8454.393199999999
$ws.Range("M134").Value = -4821.219599999999
$ws.Range("H134").Value = 2344.06
$ws.Range("I134").Value = 2452.0732
$ws.Range("K134").Value = 7356.219599999999
$ws.Range("L136").Value = 8107.799999999999
$ws.Range("M136").Value = -2278.9998
$ws.Range("H136").Value = 1798.1034
$ws.Range("I136").Value = 1609.6666
$ws.Range("K136").Value = 4828.9998
$ws.Range("J136").Value = 2702.6
$ws.Range("N136").Value = -13207.8
$ws.Range("J139").Value = 77777
$ws.Range("N139").Value = -88057
$ws.Range("L139").Value = 77777
$ws.Range("H139").Value = 77777
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M70").Value = -11661.9999
$ws.Range("H70").Value = 4959.5557
$ws.Range("I70").Value = 3992.3333
$ws.Range("K70").Value = 11976.9999
$ws.Range("H73").Value = 4959.5557
$ws.Range("I73").Value = 3992.3333
$ws.Range("K73").Value = 11976.9999
$ws.Range("M73").Value = -10884.9999
$ws.Range("N75").Value = -13003.9
$ws.Range("L75").Value = 11007.9
$ws.Range("M75").Value = -2002
$ws.Range("H75").Value = 3224.4167
$ws.Range("I75").Value = 1000
$ws.Range("K75").Value = 3000
$ws.Range("J75").Value = 3669.3
$ws.Range("K78").Value = 9000
$ws.Range("J78").Value = 3669.3
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -43007.7
$ws.Range("L78").Value = 33023.7
$ws.Range("H78").Value = 3224.4167
$ws.Range("I78").Value = 1000
$ws.Range("J86").Value = 3100
$ws.Range("N86").Value = -11672
$ws.Range("L86").Value = 9300
$ws.Range("H86").Value = 3100
$ws.Range("L89").Value = 27900
$ws.Range("H89").Value = 3100
$ws.Range("J89").Value = 3100
$ws.Range("N89").Value = -39756
$ws.Range("L92").Value = 1055.25
$ws.Range("M92").Value = -696.75
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 648.25
$ws.Range("K92").Value = 1944.75
$ws.Range("J92").Value = 351.75
$ws.Range("N92").Value = -3551.25
$ws.Range("N107").Value = -5816.470649999999
$ws.Range("L107").Value = 1976.47065
$ws.Range("M107").Value = 393
$ws.Range("H107").Value = 639.61536
$ws.Range("I107").Value = 509
$ws.Range("J107").Value = 658.82355
$ws.Range("K107").Value = 1527
$ws.Range("N122").Value = -8027.5
$ws.Range("L122").Value = 3127.5
$ws.Range("H122").Value = 503
$ws.Range("J122").Value = 347.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 29429
$ws.Range("J49").Value = 29429
$ws.Range("N49").Value = -29797
$ws.Range("L49").Value = 29429
$ws.Range("N102").Value = -4760.6666
$ws.Range("L102").Value = 1516.6666
$ws.Range("M102").Value = -99.78580000000011
$ws.Range("H102").Value = 1685.5883
$ws.Range("I102").Value = 1721.7858
$ws.Range("K102").Value = 1721.7858
$ws.Range("J102").Value = 1516.6666
$ws.Range("N122").Value = -14046.625
$ws.Range("L122").Value = 9146.625
$ws.Range("H122").Value = 2471.3044
$ws.Range("J122").Value = 3048.875
$ws.Range("N124").Value = -43820
$ws.Range("L124").Value = 34000
$ws.Range("H124").Value = 34000
$ws.Range("J124").Value = 34000
$ws.Range("I126").Value = 2666.5
$ws.Range("K126").Value = 7999.5
$ws.Range("M126").Value = -5529.5
$ws.Range("H126").Value = 2984.7778
$ws.Range("H132").Value = 1848.5294
$ws.Range("I132").Value = 1895.625
$ws.Range("K132").Value = 5686.875
$ws.Range("M132").Value = -3156.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 1837.5
$ws.Range("K22").Value = 897.75
$ws.Range("N22").Value = -2427.5
$ws.Range("M22").Value = -602.75
$ws.Range("L22").Value = 1837.5
$ws.Range("H22").Value = 1273.65
$ws.Range("I22").Value = 897.75
$ws.Range("K27").Value = 897.75
$ws.Range("J27").Value = 1837.5
$ws.Range("M27").Value = -790.75
$ws.Range("N27").Value = -2051.5
$ws.Range("L27").Value = 1837.5
$ws.Range("I27").Value = 897.75
$ws.Range("H27").Value = 1273.65
$ws.Range("K68").Value = 1516
$ws.Range("M68").Value = -767
$ws.Range("H68").Value = 1949.5714
$ws.Range("I68").Value = 1516
$ws.Range("K71").Value = 7580
$ws.Range("M71").Value = -3836
$ws.Range("H71").Value = 1949.5714
$ws.Range("I71").Value = 1516
$ws.Range("I100").Value = 1365.8889
$ws.Range("K100").Value = 1365.8889
$ws.Range("J100").Value = 52975.715
$ws.Range("M100").Value = -824.8888999999999
$ws.Range("N100").Value = -54057.715
$ws.Range("L100").Value = 52975.715
$ws.Range("H100").Value = 23945.188
$ws.Range("H132").Value = 3457.0513
$ws.Range("I132").Value = 3146.5186
$ws.Range("K132").Value = 9439.5558
$ws.Range("M132").Value = -6909.5558
$ws.Range("M136").Value = -7455.714
$ws.Range("H136").Value = 3756.9033
$ws.Range("I136").Value = 3335.238
$ws.Range("K136").Value = 10005.714
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I100").Value = 3323
$ws.Range("K100").Value = 6646
$ws.Range("J100").Value = 1629.5
$ws.Range("M100").Value = -6105
$ws.Range("N100").Value = -4341
$ws.Range("L100").Value = 3259
$ws.Range("H100").Value = 2758.5
$ws.Range("H132").Value = 1647519.8
$ws.Range("I132").Value = 7019.952
$ws.Range("K132").Value = 21059.856
$ws.Range("M132").Value = -18529.856
$ws.Range("L136").Value = 11158.8
$ws.Range("M136").Value = -447
$ws.Range("H136").Value = 3266.1667
$ws.Range("I136").Value = 999
$ws.Range("K136").Value = 2997
$ws.Range("J136").Value = 3719.6
$ws.Range("N136").Value = -16258.8
$ws.Range("J139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("L139").Value = 0
$ws.Range("H139").Value = 0